# Generate Report for Handback
# This script reflects that the de-de (and zh-cn) localization handback has
# completed and is now in sync with en-US: it
#   * updates the overall Status text ("Ready for handoff" -> "Handed back: in sync with en-US")
#   * fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#     columns on the zh-cn and de-de detail sheets (and links them back to the source doc)
#   * widens a few columns so the new, longer text fits

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94ea3feccb7950fe4fac33b8942e840c036629ee/e2e/e3b03dc6-0ded-40e2-9dc8-0e2836428ce9.md"
$mdDisplay = "e3b03dc6-0ded-40e2-9dc8-0e2836428ce9.md"

# ---------------------------------------------------------------------------
# Overview sheet: status text + wider zh-cn / de-de columns
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column now reflects the handback
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Target File / Latest Handback File / Latest Handback DateTime
$wsZh.Range("J2").Value = "e3b03dc6-0ded-40e2-9dc8-0e2836428ce9.86947f205a4d612b10d38ff70ddd3f99af248953.zh-cn.xlf"
$wsZh.Range("J3").Value = "e3b03dc6-0ded-40e2-9dc8-0e2836428ce9.86947f205a4d612b10d38ff70ddd3f99af248953.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-18 12:47:18"
$wsZh.Range("K3").Value = "2016-10-18 12:47:18"

# Widen columns to fit new content
$wsZh.Columns.Item(3).ColumnWidth  = 29.1
$wsZh.Columns.Item(9).ColumnWidth  = 39.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

# Rebuild hyperlinks so the final order/ids read A2, I2, A3, I3
$hlA2Address = $wsZh.Hyperlinks.Item(1).Address
$hlA2Display = $wsZh.Hyperlinks.Item(1).TextToDisplay
$hlA3Address = $wsZh.Hyperlinks.Item(2).Address
$hlA3Display = $wsZh.Hyperlinks.Item(2).TextToDisplay
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hlA2Address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $hlA2Display) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hlA3Address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $hlA3Display) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay) | Out-Null

# I column now holds the linked md file name too
$wsZh.Range("I2").Value = $mdDisplay
$wsZh.Range("I3").Value = $mdDisplay

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column now reflects the handback
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Latest Target File / Latest Handback File / Latest Handback DateTime
$wsDe.Range("J2").Value = "e3b03dc6-0ded-40e2-9dc8-0e2836428ce9.86947f205a4d612b10d38ff70ddd3f99af248953.de-de.xlf"
$wsDe.Range("J3").Value = "e3b03dc6-0ded-40e2-9dc8-0e2836428ce9.86947f205a4d612b10d38ff70ddd3f99af248953.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-18 12:47:35"
$wsDe.Range("K3").Value = "2016-10-18 12:47:35"

# Widen columns to fit new content
$wsDe.Columns.Item(3).ColumnWidth  = 29.1
$wsDe.Columns.Item(9).ColumnWidth  = 39.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15

# Rebuild hyperlinks so the final order/ids read A2, I2, A3, I3
$hlA2AddressDe = $wsDe.Hyperlinks.Item(1).Address
$hlA2DisplayDe = $wsDe.Hyperlinks.Item(1).TextToDisplay
$hlA3AddressDe = $wsDe.Hyperlinks.Item(2).Address
$hlA3DisplayDe = $wsDe.Hyperlinks.Item(2).TextToDisplay
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hlA2AddressDe, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $hlA2DisplayDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hlA3AddressDe, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $hlA3DisplayDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay) | Out-Null

# I column now holds the linked md file name too
$wsDe.Range("I2").Value = $mdDisplay
$wsDe.Range("I3").Value = $mdDisplay
